$d = $word.ActiveDocument

$replacements = @(
    @("2023-10-10 Tuesday", "2023-10-11 Wednesday"),
    @("27×79=2133", "58×89=5162"),
    @("77×53=4081", "57×66=3762"),
    @("61×93=5673", "44×36=1584"),
    @("62×66=4092", "89×55=4895"),
    @("23×65=1495", "50×77=3850"),
    @("26×33=858", "79×18=1422"),
    @("63×60=3780", "53×48=2544"),
    @("15×41=615", "80×51=4080"),
    @("71×18=1278", "42×15=630"),
    @("13×74=962", "60×22=1320"),
    @("69×83=5727", "72×50=3600"),
    @("92×34=3128", "84×11=924"),
    @("89×41=3649", "53×55=2915"),
    @("98×48=4704", "98×90=8820"),
    @("42×67=2814", "93×40=3720"),
    @("12×60=720", "19×67=1273"),
    @("83×97=8051", "92×13=1196"),
    @("87×47=4089", "84×56=4704"),
    @("28×80=2240", "24×49=1176"),
    @("77×22=1694", "69×69=4761"),
    @("76×66=5016", "13×72=936"),
    @("24×30=720", "63×41=2583"),
    @("72×41=2952", "71×55=3905"),
    @("67×20=1340", "97×40=3880"),
    @("13×61=793", "66×27=1782")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
